$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 2.25
$ws.Range("H3").Value = 2.9
$ws.Range("I3").Value = 3.6
$ws.Range("K3").Value = 1.91
$ws.Range("L3").Value = 4.33
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5
$ws.Range("O3").Value = 1.53
$ws.Range("P3").Value = 2.38
$ws.Range("Q3").Value = 2.7
$ws.Range("R3").Value = 1.44
$ws.Range("AC3").Value = 6
$ws.Range("AI3").Value = 13
$ws.Range("AK3").Value = 34
$ws.Range("AO3").Value = 13
$ws.Range("AP3").Value = 29
$ws.Range("AS3").Value = 301
$ws.Range("AU3").Value = 9.5
$ws.Range("AX3").Value = 21
$ws.Range("BA3").Value = 126
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 3.4
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 7
$ws.Range("O4").Value = 1.4
$ws.Range("P4").Value = 2.75
$ws.Range("AC4").Value = 7
$ws.Range("AF4").Value = 51
$ws.Range("AG4").Value = 8.5
$ws.Range("AI4").Value = 13
$ws.Range("AM4").Value = 1000
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 9
$ws.Range("O5").Value = 1.36
$ws.Range("P5").Value = 3
$ws.Range("Q5").Value = 2.2
$ws.Range("R5").Value = 1.65
$ws.Range("G8").Value = 2.05
$ws.Range("H8").Value = 3.7
$ws.Range("I8").Value = 3.25
$ws.Range("O8").Value = 1.18
$ws.Range("P8").Value = 4.5
$ws.Range("Q8").Value = 1.6
$ws.Range("R8").Value = 2.3
$ws.Range("U8").Value = 1.5
$ws.Range("V8").Value = 2.5
$ws.Range("Z8").Value = 19
$ws.Range("AC8").Value = 17
$ws.Range("AM8").Value = 101
$ws.Range("AN8").Value = 4.33
$ws.Range("AR8").Value = 41
$ws.Range("M11").Value = 1.06
$ws.Range("N11").Value = 10
$ws.Range("O11").Value = 1.33
$ws.Range("P11").Value = 3.25
$ws.Range("H13").Value = 3.3
$ws.Range("I13").Value = 6.5
$ws.Range("L13").Value = 7.5
$ws.Range("AI13").Value = 23
$ws.Range("AU13").Value = 12
$ws.Range("H14").Value = 3.9
$ws.Range("I14").Value = 4
$ws.Range("N14").Value = 17
$ws.Range("O14").Value = 1.18
$ws.Range("P14").Value = 4.5
$ws.Range("AC14").Value = 17
$ws.Range("AD14").Value = 8
$ws.Range("AE14").Value = 13
$ws.Range("AG14").Value = 15
$ws.Range("AI14").Value = 13
$ws.Range("AL14").Value = 29
$ws.Range("H15").Value = 2.9
$ws.Range("M15").Value = 1.11
$ws.Range("N15").Value = 6.5
$ws.Range("O15").Value = 1.53
$ws.Range("P15").Value = 2.38
$ws.Range("Q15").Value = 2.7
$ws.Range("R15").Value = 1.44
$ws.Range("AA15").Value = 29
$ws.Range("AR15").Value = 101
$ws.Range("M16").Value = 1.11
$ws.Range("N16").Value = 6.5

Write-Host "Applied odds updates"
